# Auto-generated edit script applying targeted cell value updates
# across multiple sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as described by the commit's underlying data refresh.

$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1029.6364
$ws.Range("I19").Value = 908.5
$ws.Range("K19").Value = 908.5
$ws.Range("M19").Value = -733.5
$ws.Range("H28").Value = 1017.63635
$ws.Range("I28").Value = 954.1111
$ws.Range("K28").Value = 954.1111
$ws.Range("M28").Value = -469.1111
$ws.Range("H40").Value = 6900
$ws.Range("I40").Value = 4075
$ws.Range("K40").Value = 4075
$ws.Range("M40").Value = -3900
$ws.Range("H55").Value = 190
$ws.Range("I55").Value = 190
$ws.Range("K55").Value = 190
$ws.Range("M55").Value = 24
$ws.Range("H58").Value = 1404.75
$ws.Range("J58").Value = 2999
$ws.Range("L58").Value = 8997
$ws.Range("N58").Value = -9297
$ws.Range("H64").Value = 4234.857
$ws.Range("I64").Value = 4419.5
$ws.Range("J64").Value = 3773.25
$ws.Range("K64").Value = 4419.5
$ws.Range("L64").Value = 3773.25
$ws.Range("M64").Value = -4171.5
$ws.Range("N64").Value = -4269.25
$ws.Range("H67").Value = 4234.857
$ws.Range("I67").Value = 4419.5
$ws.Range("J67").Value = 3773.25
$ws.Range("K67").Value = 4419.5
$ws.Range("L67").Value = 3773.25
$ws.Range("M67").Value = -3561.5
$ws.Range("N67").Value = -5489.25
$ws.Range("H74").Value = 100006320
$ws.Range("I74").Value = 166674370
$ws.Range("K74").Value = 166674370
$ws.Range("M74").Value = -166673434
$ws.Range("H77").Value = 100006320
$ws.Range("I77").Value = 166674370
$ws.Range("K77").Value = 833371850
$ws.Range("M77").Value = -833367170
$ws.Range("H88").Value = 12088.444
$ws.Range("I88").Value = 8999
$ws.Range("J88").Value = 12971.143
$ws.Range("K88").Value = 8999
$ws.Range("L88").Value = 12971.143
$ws.Range("M88").Value = -8593
$ws.Range("N88").Value = -13783.143
$ws.Range("H91").Value = 12088.444
$ws.Range("I91").Value = 8999
$ws.Range("J91").Value = 12971.143
$ws.Range("K91").Value = 8999
$ws.Range("L91").Value = 12971.143
$ws.Range("M91").Value = -7595
$ws.Range("N91").Value = -15779.143
$ws.Range("H138").Value = 2774.7017
$ws.Range("J138").Value = 2937.1956
$ws.Range("L138").Value = 8811.586800000001
$ws.Range("N138").Value = -19091.5868

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2815.1667
$ws.Range("J2").Value = 2736.75
$ws.Range("L2").Value = 2736.75
$ws.Range("N2").Value = -2962.75
$ws.Range("H32").Value = 3616.182
$ws.Range("I32").Value = 2359.8333
$ws.Range("K32").Value = 2359.8333
$ws.Range("M32").Value = -2072.8333
$ws.Range("H45").Value = 4710.391
$ws.Range("I45").Value = 6571.769
$ws.Range("K45").Value = 6571.769
$ws.Range("M45").Value = -6194.769
$ws.Range("H74").Value = 22731440
$ws.Range("I74").Value = 25003960
$ws.Range("K74").Value = 25003960
$ws.Range("M74").Value = -25003086
$ws.Range("H77").Value = 22731440
$ws.Range("I77").Value = 25003960
$ws.Range("K77").Value = 125019800
$ws.Range("M77").Value = -125015432
$ws.Range("H88").Value = 3372.75
$ws.Range("I88").Value = 2998.8
$ws.Range("J88").Value = 3996
$ws.Range("K88").Value = 2998.8
$ws.Range("L88").Value = 3996
$ws.Range("M88").Value = -2592.8
$ws.Range("N88").Value = -4808
$ws.Range("H91").Value = 3372.75
$ws.Range("I91").Value = 2998.8
$ws.Range("J91").Value = 3996
$ws.Range("K91").Value = 2998.8
$ws.Range("L91").Value = 3996
$ws.Range("M91").Value = -1594.8
$ws.Range("N91").Value = -6804
$ws.Range("H116").Value = 2815.1667
$ws.Range("J116").Value = 2736.75
$ws.Range("L116").Value = 2736.75
$ws.Range("N116").Value = -7324.75
$ws.Range("H122").Value = 7803.385
$ws.Range("I122").Value = 7777.7144
$ws.Range("K122").Value = 23333.1432
$ws.Range("M122").Value = -20883.1432
$ws.Range("H132").Value = 5887309.5
$ws.Range("I132").Value = 6254641
$ws.Range("K132").Value = 18763923
$ws.Range("M132").Value = -18761393

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2815.1667
$ws.Range("J3").Value = 2736.75
$ws.Range("L3").Value = 2736.75
$ws.Range("N3").Value = -2964.75
$ws.Range("H134").Value = 19234208
$ws.Range("I134").Value = 20836624
$ws.Range("K134").Value = 62509872
$ws.Range("M134").Value = -62507337

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H31").Value = 2791
$ws.Range("I31").Value = 2715.2666
$ws.Range("K31").Value = 2715.2666
$ws.Range("M31").Value = -2420.2666
$ws.Range("H34").Value = 2791
$ws.Range("I34").Value = 2715.2666
$ws.Range("K34").Value = 2715.2666
$ws.Range("M34").Value = -2513.2666
$ws.Range("H105").Value = 3242
$ws.Range("I105").Value = 3905
$ws.Range("J105").Value = 2800
$ws.Range("K105").Value = 3905
$ws.Range("L105").Value = 2800
$ws.Range("M105").Value = -2158
$ws.Range("N105").Value = -6294
$ws.Range("H132").Value = 52633896
$ws.Range("I132").Value = 66669120
$ws.Range("K132").Value = 200007360
$ws.Range("M132").Value = -200004830
$ws.Range("H134").Value = 31381106
$ws.Range("I134").Value = 35863010
$ws.Range("K134").Value = 107589030
$ws.Range("M134").Value = -107586495

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 14583.333
$ws.Range("I57").Value = 2500
$ws.Range("J57").Value = 17000
$ws.Range("K57").Value = 7500
$ws.Range("L57").Value = 51000
$ws.Range("M57").Value = -6941
$ws.Range("N57").Value = -52118
$ws.Range("H80").Value = 3363.5
$ws.Range("I80").Value = 3596.4
$ws.Range("K80").Value = 10789.2
$ws.Range("M80").Value = -9853.200000000001
$ws.Range("H83").Value = 3363.5
$ws.Range("I83").Value = 3596.4
$ws.Range("K83").Value = 32367.6
$ws.Range("M83").Value = -27687.6
$ws.Range("H131").Value = 1636.3556
$ws.Range("I131").Value = 890.375
$ws.Range("K131").Value = 2671.125
$ws.Range("M131").Value = 2368.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 9622242
$ws.Range("I132").Value = 11369931
$ws.Range("K132").Value = 34109793
$ws.Range("M132").Value = -34107263

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2943.625
$ws.Range("I22").Value = 3591.6667
$ws.Range("J22").Value = 999.5
$ws.Range("K22").Value = 3591.6667
$ws.Range("L22").Value = 999.5
$ws.Range("M22").Value = -3296.6667
$ws.Range("N22").Value = -1589.5
$ws.Range("H27").Value = 2943.625
$ws.Range("I27").Value = 3591.6667
$ws.Range("J27").Value = 999.5
$ws.Range("K27").Value = 3591.6667
$ws.Range("L27").Value = 999.5
$ws.Range("M27").Value = -3484.6667
$ws.Range("N27").Value = -1213.5
$ws.Range("H46").Value = 1309
$ws.Range("I46").Value = 1336.5
$ws.Range("K46").Value = 1336.5
$ws.Range("M46").Value = -1148.5
$ws.Range("H68").Value = 2634027.2
$ws.Range("I68").Value = 3761725.5
$ws.Range("K68").Value = 3761725.5
$ws.Range("M68").Value = -3760976.5
$ws.Range("H71").Value = 2634027.2
$ws.Range("I71").Value = 3761725.5
$ws.Range("K71").Value = 18808627.5
$ws.Range("M71").Value = -18804883.5
$ws.Range("H132").Value = 125005000
$ws.Range("I132").Value = 125005000
$ws.Range("K132").Value = 375015000
$ws.Range("M132").Value = -375012470

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1961.5555
$ws.Range("I81").Value = 1093.5714
$ws.Range("K81").Value = 2187.1428
$ws.Range("M81").Value = -1126.1428
$ws.Range("H84").Value = 1961.5555
$ws.Range("I84").Value = 1093.5714
$ws.Range("K84").Value = 10935.714
$ws.Range("M84").Value = -5631.714
$ws.Range("H132").Value = 13520603
$ws.Range("I132").Value = 18521620
$ws.Range("K132").Value = 55564860
$ws.Range("M132").Value = -55562330

Write-Host "Applied 206 cell updates across 8 sheets"
